$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8 ("date"), pushing the rest down by one.
$ws.Rows.Item(8).Insert()

# Copy formatting (styles) from the row that is now row 9 (the former row 8, "date")
# onto the newly inserted blank row 8, so the new row matches the sheet's normal styling.
$ws.Range("B9:E9").Copy()
$ws.Range("B8:E8").PasteSpecial(-4122)

# Populate the new row with the "numeric" type and its restriction text.
$ws.Range("B8").Value = "numeric"
$ws.Range("C8").Value = "`"type`": [ `"number`", `"null`" ],`n`"minimum`": 1"

# Match the row height used for this new entry.
$ws.Rows.Item(8).RowHeight = 25.5

# Update the selected/active cell to reflect the new layout (row shifted from 10 to 9).
[void]$ws.Range("C9").Select()
